# Weekly update: insert a new data row at row 70 (pushing existing rows
# 70-174 down to 71-175) and populate the new row with this week's
# "Repollo" (Crespo record) price observation for Terminal
# Hortofrutícola Agro Chillán.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 70; this shifts rows 70-174
# down to 71-175, preserving all of their data/formatting untouched.
$ws.Rows.Item(70).Insert()

# Populate the newly-inserted row 70 with the new weekly observation.
$ws.Range("A70").Value = 7
$ws.Range("B70").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C70").Value = "Ñuble"
$ws.Range("D70").Value = 44579
$ws.Range("E70").Value = 16
$ws.Range("F70").Value = 100112006
$ws.Range("G70").Value = "Repollo"
$ws.Range("H70").Value = "Crespo record"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 300
$ws.Range("K70").Value = 600
$ws.Range("L70").Value = 700
$ws.Range("M70").Value = 650
$ws.Range("N70").Value = "$/unidad"
$ws.Range("O70").Value = "Provincia de Diguillín"
$ws.Range("P70").Value = 650
$ws.Range("Q70").Value = 1
$ws.Range("R70").Value = "Hortaliza"
